# Apply predictor-label edits described in the diff:
# wrap several predictor names with "ln(...)" (using square brackets inside),
# and fix the "Livestock AB Consumption" label bracket typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new predictor text (column C holds "Predictor" values)
$ws.Range("C2").Value  = "ln(GDP [dollars per capita])"
$ws.Range("C17").Value = "ln(GDP [dollars per capita])"

$ws.Range("C3").Value  = "ln(Migrant Population [per capita])"

$ws.Range("C4").Value  = "ln(Tourism - Inbound [per capita])"

$ws.Range("C5").Value  = "ln(ProMed Mentions [per capita])"
$ws.Range("C13").Value = "ln(ProMed Mentions [per capita])"

$ws.Range("C8").Value  = "ln(Publication Bias Index [per capita])"
$ws.Range("C15").Value = "ln(Publication Bias Index [per capita])"

$ws.Range("C9").Value  = "ln(AB Exports [dollars per capita])"

$ws.Range("C12").Value = "Livestock AB Consumption [kg per capita)"

$ws.Range("C16").Value = "ln(Population)"

$wb.Save()
